$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 65000
$ws.Range("J3").Value = 65000
$ws.Range("L3").Value = 65000
$ws.Range("N3").Value = -65228
$ws.Range("H19").Value = 986.5454999999999
$ws.Range("I19").Value = 1023.1429
$ws.Range("K19").Value = 1023.1429
$ws.Range("M19").Value = -848.1429000000001
$ws.Range("H55").Value = 500
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 500
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 500
$ws.Range("M55").Value = ""
$ws.Range("N55").Value = -928
$ws.Range("H96").Value = 526.2273
$ws.Range("I96").Value = 309.4737
$ws.Range("J96").Value = 1899
$ws.Range("K96").Value = 928.4211
$ws.Range("L96").Value = 5697
$ws.Range("M96").Value = 444.5789
$ws.Range("N96").Value = -8443
$ws.Range("H100").Value = 41449.58
$ws.Range("I100").Value = 46412.652
$ws.Range("J100").Value = 3399.3333
$ws.Range("K100").Value = 46412.652
$ws.Range("L100").Value = 3399.3333
$ws.Range("M100").Value = -45871.652
$ws.Range("N100").Value = -4481.3333
$ws.Range("H102").Value = 65000
$ws.Range("J102").Value = 65000
$ws.Range("L102").Value = 65000
$ws.Range("N102").Value = -71490
$ws.Range("H111").Value = 952.2308
$ws.Range("I111").Value = 952.2308
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 2856.6924
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 210.3076000000001
$ws.Range("N111").Value = ""
$ws.Range("H112").Value = 1599.9395
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1599.9395
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 4799.818499999999
$ws.Range("M112").Value = ""
$ws.Range("N112").Value = -7015.818499999999
$ws.Range("H132").Value = 1974.7142
$ws.Range("I132").Value = 1865
$ws.Range("J132").Value = 2249
$ws.Range("K132").Value = 5595
$ws.Range("L132").Value = 6747
$ws.Range("M132").Value = -3065
$ws.Range("N132").Value = -11807
$ws.Range("H135").Value = 2415
$ws.Range("I135").Value = 2415
$ws.Range("K135").Value = 21735
$ws.Range("M135").Value = -19200
$ws.Range("H138").Value = 2660.9395
$ws.Range("J138").Value = 3369.818
$ws.Range("L138").Value = 10109.454
$ws.Range("N138").Value = -20389.454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7452.7617
$ws.Range("I45").Value = 8183.25
$ws.Range("J45").Value = 5115.2
$ws.Range("K45").Value = 8183.25
$ws.Range("L45").Value = 5115.2
$ws.Range("M45").Value = -7806.25
$ws.Range("N45").Value = -5869.2
$ws.Range("H61").Value = 6448.0625
$ws.Range("I61").Value = 3346.1
$ws.Range("K61").Value = 3346.1
$ws.Range("M61").Value = -3134.1
$ws.Range("H74").Value = 3379.05
$ws.Range("I74").Value = 2661.9092
$ws.Range("J74").Value = 4255.5557
$ws.Range("K74").Value = 2661.9092
$ws.Range("L74").Value = 4255.5557
$ws.Range("M74").Value = -1787.9092
$ws.Range("N74").Value = -6003.5557
$ws.Range("H77").Value = 3379.05
$ws.Range("I77").Value = 2661.9092
$ws.Range("J77").Value = 4255.5557
$ws.Range("K77").Value = 13309.546
$ws.Range("L77").Value = 21277.7785
$ws.Range("M77").Value = -8941.546
$ws.Range("N77").Value = -30013.7785
$ws.Range("H103").Value = 36200
$ws.Range("J103").Value = 36200
$ws.Range("L103").Value = 36200
$ws.Range("N103").Value = -38544
$ws.Range("H110").Value = 1908.6207
$ws.Range("I110").Value = 1826.04
$ws.Range("J110").Value = 2424.75
$ws.Range("K110").Value = 1826.04
$ws.Range("L110").Value = 2424.75
$ws.Range("M110").Value = 218.96
$ws.Range("N110").Value = -6514.75
$ws.Range("H132").Value = 5717.5435
$ws.Range("I132").Value = 5809.067
$ws.Range("K132").Value = 17427.201
$ws.Range("M132").Value = -14897.201
$ws.Range("H136").Value = 6448.0625
$ws.Range("I136").Value = 3346.1
$ws.Range("K136").Value = 10038.3
$ws.Range("M136").Value = -7488.299999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = -127
$ws.Range("H92").Value = 63249.25
$ws.Range("J92").Value = 63249.25
$ws.Range("L92").Value = 63249.25
$ws.Range("N92").Value = -68241.25
$ws.Range("H122").Value = 30780
$ws.Range("J122").Value = 30780
$ws.Range("L122").Value = 30780
$ws.Range("N122").Value = -40580

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3059.7568
$ws.Range("I16").Value = 3111.8147
$ws.Range("K16").Value = 3111.8147
$ws.Range("M16").Value = -2824.8147
$ws.Range("H43").Value = 11921.2
$ws.Range("J43").Value = 11921.2
$ws.Range("L43").Value = 11921.2
$ws.Range("N43").Value = -12289.2
$ws.Range("H95").Value = 38539.6
$ws.Range("J95").Value = 38539.6
$ws.Range("L95").Value = 38539.6
$ws.Range("N95").Value = -44031.6
$ws.Range("H101").Value = 11921.2
$ws.Range("J101").Value = 11921.2
$ws.Range("L101").Value = 11921.2
$ws.Range("N101").Value = -18411.2
$ws.Range("H107").Value = 728.8333
$ws.Range("I107").Value = 624.4286
$ws.Range("J107").Value = 1459.6666
$ws.Range("K107").Value = 624.4286
$ws.Range("L107").Value = 1459.6666
$ws.Range("M107").Value = 1295.5714
$ws.Range("N107").Value = -5299.6666
$ws.Range("H113").Value = 3059.7568
$ws.Range("I113").Value = 3111.8147
$ws.Range("K113").Value = 3111.8147
$ws.Range("M113").Value = -941.8146999999999
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").Value = ""
$ws.Range("H132").Value = 22910.984
$ws.Range("I132").Value = 14010.652
$ws.Range("K132").Value = 42031.956
$ws.Range("M132").Value = -39501.956
$ws.Range("H140").Value = 94128.664
$ws.Range("J140").Value = 99995
$ws.Range("L140").Value = 99995
$ws.Range("N140").Value = -110355

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2713.7354
$ws.Range("I34").Value = 116.42857
$ws.Range("K34").Value = 349.28571
$ws.Range("M34").Value = -265.28571
$ws.Range("H131").Value = 16588.924
$ws.Range("J131").Value = 1482.3281
$ws.Range("L131").Value = 4446.9843
$ws.Range("N131").Value = -14526.9843

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1819.8
$ws.Range("J80").Value = 2239.8
$ws.Range("L80").Value = 2239.8
$ws.Range("N80").Value = -4235.8
$ws.Range("H83").Value = 1819.8
$ws.Range("J83").Value = 2239.8
$ws.Range("L83").Value = 11199
$ws.Range("N83").Value = -21183
$ws.Range("H132").Value = 2861.2222
$ws.Range("I132").Value = 2861.2222
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8583.6666
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6053.6666
$ws.Range("N132").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2099.7
$ws.Range("I61").Value = 2099.7
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2099.7
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1897.7
$ws.Range("N61").Value = ""
$ws.Range("H104").Value = 15294.833
$ws.Range("J104").Value = 15294.833
$ws.Range("L104").Value = 15294.833
$ws.Range("N104").Value = -22282.833
$ws.Range("H113").Value = 2099.7
$ws.Range("I113").Value = 2099.7
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2099.7
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 70.30000000000018
$ws.Range("N113").Value = ""
$ws.Range("H132").Value = 6539491
$ws.Range("I132").Value = 6947646.5
$ws.Range("K132").Value = 20842939.5
$ws.Range("M132").Value = -20840409.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1580.963
$ws.Range("I107").Value = 1653.5385
$ws.Range("J107").Value = 1513.5714
$ws.Range("K107").Value = 4960.6155
$ws.Range("L107").Value = 4540.7142
$ws.Range("M107").Value = -3040.6155
$ws.Range("N107").Value = -8380.7142
$ws.Range("H113").Value = 600.4
$ws.Range("I113").Value = 528.55554
$ws.Range("K113").Value = 1585.66662
$ws.Range("M113").Value = 584.33338
$ws.Range("H132").Value = 19322.096
$ws.Range("I132").Value = 12975.892
$ws.Range("K132").Value = 38927.676
$ws.Range("M132").Value = -36397.676
